$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.698.65"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.531.41"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'205.46"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "'0.484"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'21.35"
$ws.Range("E8").Value = "  -2.59%  "
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").Value = "1.748.75"
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("D13").Value = "1.534.55"
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").Value = "'3.67"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").Value = "'0.505"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").Value = "26.683.39"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "'61.04"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "'212.16"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").Value = "0.0₃0681"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  -2.22%  "
$ws.Range("E23").Value = "  -3.18%  "
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("D25").Value = "'151.73"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").Value = "'6.57"
$ws.Range("E26").Value = "  -2.92%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("D31").Value = "'0.0452"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").Value = "'3.23"
$ws.Range("E32").Value = "  +2.72%  "
$ws.Range("D33").Value = "1.358.85"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").Value = "'1.50"
$ws.Range("E35").Value = "  -2.98%  "
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "'0.936"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'0.797"
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("D42").Value = "'5.69"
$ws.Range("E42").Value = "  +5.94%  "
$ws.Range("D43").Value = "'0.990"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D45").Value = "'1.74"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("D46").Value = "'62.35"
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").Value = "1.663.13"
$ws.Range("E47").Value = "  -1.78%  "
$ws.Range("D48").Value = "'85.31"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("D50").Value = "0.0₇0970"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").Value = "'0.0941"
$ws.Range("E51").Value = "  -0.70%  "
